$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the scraped crypto market update for this run.
# Columns B (Coin), C (Link), D (Price) and E (Volume 1h) are updated
# per-row. D/B/C are forced to Text format before assignment so that
# numeric-looking strings (e.g. "1.00", "13.00", "0.0000240") are not
# auto-coerced into numbers/scientific notation by Excel, preserving
# the exact text content.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value2 = '64.078.94'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value2 = '  +0.22%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value2 = '3.089.41'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value2 = '  +0.82%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value2 = '  -0.77%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '595.71'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value2 = '  +2.18%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value2 = '156.25'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value2 = '  +2.76%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value2 = '1.00'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value2 = '  -0.42%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value2 = '0.540'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value2 = '  +2.05%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value2 = '3.090.09'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value2 = '  +1.01%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value2 = '0.158'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value2 = '  +1.32%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value2 = '5.93'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value2 = '  +1.15%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value2 = '0.454'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value2 = '  +0.07%  '

$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value2 = 'Avalanche'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value2 = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value2 = '37.18'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value2 = '  +0.07%  '

$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value2 = 'ShibaInu'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value2 = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value2 = '0.0000240'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value2 = '  -0.17%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value2 = '0.121'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value2 = '  +1.72%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value2 = '3.609.01'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value2 = '  +0.60%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value2 = '7.24'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value2 = '  +2.11%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value2 = '64.038.85'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value2 = '  +0.41%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '3.088.85'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value2 = '  -1.37%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value2 = '489.19'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value2 = '  +5.53%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value2 = '14.66'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value2 = '  +1.21%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '0.714'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value2 = '  -1.10%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value2 = '7.65'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value2 = '  +2.56%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value2 = '2.45'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value2 = '  +4.38%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value2 = '82.21'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value2 = '  +1.63%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value2 = '13.00'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value2 = '  -0.84%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value2 = '10.68'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value2 = '  +8.85%  '

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value2 = '  +0.22%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value2 = '7.61'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value2 = '  +4.62%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value2 = '2.27'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value2 = '  +4.60%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value2 = '2.71'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value2 = '  +1.46%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value2 = '0.998'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value2 = '  -1.07%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value2 = '27.44'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value2 = '  +1.38%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value2 = '0.113'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value2 = '  -1.01%  '

$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value2 = 'PEPE'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value2 = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value2 = '0.0₃0832'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value2 = '  -0.89%  '

$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value2 = 'Mantle'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value2 = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value2 = '1.07'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value2 = '  +2.10%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value2 = '6.11'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value2 = '  +0.87%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value2 = '2.26'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value2 = '  +1.66%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value2 = '3.27'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value2 = '  -2.13%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value2 = '9.33'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value2 = '  +1.50%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value2 = '50.84'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value2 = '  +0.92%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value2 = '444.00'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value2 = '  -0.05%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value2 = '0.293'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value2 = '  +3.05%  '

$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value2 = 'VeChain'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value2 = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value2 = '0.0368'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value2 = '  +1.44%  '

$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value2 = 'Kaspa'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value2 = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '0.113'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value2 = '  +4.94%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value2 = '2.848.84'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value2 = '  +1.19%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value2 = '39.98'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value2 = '  +1.43%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value2 = '132.70'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value2 = '  +2.54%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value2 = '25.70'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value2 = '  +2.57%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value2 = '2.26'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value2 = '  +1.75%  '
